# Update countries & provincias Spain
#
# The source data feed re-sorted a handful of adjacent country rows and
# refreshed some case counts for the newer timestamp (26 Jun 2020, 05:35).
# Net effect on the worksheet "Pais":
#   - A handful of country-name pairs swap rows (and in two cases, B:H data
#     was also refreshed with new figures).
#   - A few rows got simple numeric updates (no name change).
#   - The "last updated" footer text changes time from 04:18 to 05:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 14: Mexico -- refreshed D/E only ---
$ws.Range("D14").Value = 116862
$ws.Range("E14").Value = 61029

# --- Row 24: China -- refreshed D/E only ---
$ws.Range("D24").Value = 78439
$ws.Range("E24").Value = 389

# --- Rows 60/61: Azerbaiyan <-> Honduras swap places, with fresh data ---
$ws.Range("A60").Value = "Honduras"
$ws.Range("B60").Value = 15366
$ws.Range("C60").Value = 795
$ws.Range("D60").Value = 1600
$ws.Range("E60").Value = 13340
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 9
$ws.Range("H60").Value = 426

$ws.Range("A61").Value = "Azerbaiyan"
$ws.Range("B61").Value = 14852
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 8059
$ws.Range("E61").Value = 6613
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 180

# --- Rows 88/89: Bulgaria <-> Venezuela swap places, with fresh data ---
$ws.Range("A88").Value = "Venezuela"
$ws.Range("B88").Value = 4563
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 1327
$ws.Range("E88").Value = 3197
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 39

$ws.Range("A89").Value = "Bulgaria"
$ws.Range("B89").Value = 4408
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 2370
$ws.Range("E89").Value = 1827
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 211

# --- Row 165: Mongolia -- refreshed B/C/D/E only ---
$ws.Range("B165").Value = 219
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 170
$ws.Range("E165").Value = 49

# --- Rows 202/203: Fiyi <-> Dominica swap places (data unchanged) ---
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# --- Rows 208/209: Groenlandia <-> Islas Malvinas swap places (data unchanged) ---
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"

# --- Rows 211/212: Seychelles <-> Montserrat swap places, with fresh data ---
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

$ws.Range("A212").Value = "Seychelles"
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

# --- Footer timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 05:35"
